# Towns workbook update:
#  - Update Madrid's lat/lng (row 3) to the new marker-icon derived coordinates.
#  - Move the active cell selection from G8 to G6.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Madrid latitude / longitude (row 3: E = lat, F = lng)
$ws.Range("E3").Value = 40.416951400000002
$ws.Range("F3").Value = -3.7057172

# Update the selected/active cell shown in the sheet view
$ws.Range("G6").Select()
